# Gantt chart update:
#  - Fill in actual Start/Due dates + weekly "X" markers for three existing
#    tasks (Light puzzle / Moving Light puzzle / Gun) that previously had
#    blank schedules.
#  - Insert two new tasks under a new "Interaction Interface" row just above
#    the "Parts and Assembly" section (rows 44-45), each with dates, percent
#    complete and weekly "X" markers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 37: "Light puzzle" ------------------------------------------------
$ws.Range("E37").Value = "11/10/2022"
$ws.Range("F37").Value = "11/13/2022"
$ws.Range("BC37:BG37").Value = "X"
$ws.Range("BR37:BZ37").Value = "X"

# ---- Row 38: "Moving Light puzzle" ----------------------------------------
$ws.Range("E38").Value = "11/15/2022"
$ws.Range("F38").Value = "11/21/2022"
$ws.Range("BH38:BK38").Value = "X"
$ws.Range("BR38:BZ38").Value = "X"

# ---- Row 39: "Gun" ----------------------------------------------------------
$ws.Range("E39").Value = "9/20/2022"
$ws.Range("F39").Value = "11/7/2022"
$ws.Range("X39:AB39").Value = "X"
$ws.Range("AC39").Value = "X"
$ws.Range("AI39:AO39").Value = "X"
$ws.Range("AW39:AZ39").Value = "X"
$ws.Range("BP39:BZ39").Value = "X"

# ---- Insert two new task rows before row 44 ("Parts and Assembly") --------
$ws.Rows.Item(44).Insert()
$ws.Rows.Item(44).Insert()

# New row 44: "Interaction Interface"
$ws.Range("C44").Value = "Interaction Interface"
$ws.Range("D44").Value = "BK"
$ws.Range("E44").Value = "9/20/2022"
$ws.Range("F44").Value = "10/1/2022"
$ws.Range("H44").Value = 100
$ws.Range("X44:AC44").Value = "X"

# New row 45: "Get Rubiks Cube with Gun Puzzle"
$ws.Range("C45").Value = "Get Rubiks Cube with Gun Puzzle"
$ws.Range("D45").Value = "BK"
$ws.Range("E45").Value = "11/22/2022"
$ws.Range("F45").Value = "12/2/2022"
$ws.Range("H45").Value = 100
$ws.Range("BM45:BQ45").Value = "X"
$ws.Range("BT45:BW45").Value = "X"
